$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 64, shifting existing rows 64..179 down to 65..180.
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new record's data.
# Columns that are constant across the whole data block (A,B,C,E,F,G,H,I,N,O,Q,R)
# are copied from the row below (which now holds the old row-64 data, itself a
# copy of the same constant values), so we just fill the full row explicitly.
$ws.Cells.Item(64, 1).Value = 10
$ws.Cells.Item(64, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(64, 3).Value = "La Araucanía"
$ws.Cells.Item(64, 4).Value = 44771
$ws.Cells.Item(64, 5).Value = 9
$ws.Cells.Item(64, 6).Value = 100114007
$ws.Cells.Item(64, 7).Value = "Jengibre"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 40
$ws.Cells.Item(64, 11).Value = 20000
$ws.Cells.Item(64, 12).Value = 20000
$ws.Cells.Item(64, 13).Value = 20000
$ws.Cells.Item(64, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(64, 15).Value = "Perú"
$ws.Cells.Item(64, 16).Value = 1538
$ws.Cells.Item(64, 17).Value = 13
$ws.Cells.Item(64, 18).Value = "Hortaliza"

# Apply the same date style (style index 2 in the original file, i.e. the
# numFmt used by column D) to the new D64 cell so it matches the rest of the
# column. Copy the number format from the cell above (D63) which already has
# the correct date format.
$ws.Cells.Item(64, 4).NumberFormat = $ws.Cells.Item(63, 4).NumberFormat
